# Assignment 1 codebook - Part C update
# Insert a new "VFA_" / "Visceral obesity" variable definition as rows 44-45,
# shifting all the subsequent variable rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two blank rows at position 44 (pushes former rows 44-55 down to 46-57)
$ws.Rows("44:45").Insert()

# Excel auto-fills the inserted rows with empty, but styled, placeholder cells
# in the columns that carry a column-level style (E:G and M). The new variable
# row doesn't use those columns, so clear them away entirely.
$ws.Range("E44:G45").Clear()
$ws.Range("M44:M45").Clear()

# Row 44: new variable "VFA_" ("Visceral obesity"), category "No"
$ws.Range("A44").Value = 28
$ws.Range("B44").Value = "VFA_"
$ws.Range("C44").Value = "Visceral obesity"
$ws.Range("D44").Value = "Numeric"
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = "No"
$ws.Range("I44").WrapText = $true

# Row 45: second category, "Yes"
$ws.Range("H45").Value = 1
$ws.Range("I45").Value = "Yes"
$ws.Range("I45").WrapText = $true

# Restore the view state (scrolled down, with the last-edited cell selected)
[void]$ws.Range("A37").Select()
[void]$ws.Range("G57").Select()
